$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text fixes (rich-text runs) ---
# A8: "Volume 30   Number  50" -> "...51"  (last run "50" -> "51")
$ws.Range("A8").Characters(21, 2).Text = "51"

# C9: "Report Covering the Week  12/11/2023  Through  12/17/2023"
#     -> "...12/18/2023  Through  12/24/2023"
$ws.Range("C9").Characters(27, 10).Text = "12/18/2023"
$ws.Range("C9").Characters(48, 10).Text = "12/24/2023"

# --- Row 30: D30/E30 switch from numeric to the same text labels used in
#     row 14 ("0" / "***.*"), including matching style. Copy D14/E14
#     (format + value) onto D30/E30 so style + shared-string text match.
$ws.Range("D14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D30").PasteSpecial(-4163)

$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4163)

$excel.CutCopyMode = 0

# --- Weekly numeric data refresh (new crime data collected) ---
$ws.Range("F14").Value = 4
$ws.Range("H14").Value = 0
$ws.Range("L14").Value = -29.032258064516
$ws.Range("M14").Value = -49.230769230769
$ws.Range("N14").Value = -85.897435897435
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 8
$ws.Range("E15").Value = -62.5
$ws.Range("G15").Value = 18
$ws.Range("H15").Value = -5.555555555555
$ws.Range("I15").Value = 213
$ws.Range("J15").Value = 250
$ws.Range("K15").Value = -14.8
$ws.Range("L15").Value = -2.293577981651
$ws.Range("M15").Value = -5.333333333333
$ws.Range("N15").Value = -63.837011884550
$ws.Range("C16").Value = 56
$ws.Range("D16").Value = 35
$ws.Range("E16").Value = 60
$ws.Range("F16").Value = 200
$ws.Range("G16").Value = 164
$ws.Range("H16").Value = 21.951219512195
$ws.Range("I16").Value = 2489
$ws.Range("J16").Value = 2503
$ws.Range("K16").Value = -0.559328805433
$ws.Range("L16").Value = 20.009643201542
$ws.Range("M16").Value = -31.938747607328
$ws.Range("N16").Value = -84.995177236556
$ws.Range("C17").Value = 75
$ws.Range("D17").Value = 76
$ws.Range("E17").Value = -1.315789473684
$ws.Range("F17").Value = 288
$ws.Range("G17").Value = 269
$ws.Range("H17").Value = 7.063197026022
$ws.Range("I17").Value = 4152
$ws.Range("J17").Value = 4070
$ws.Range("K17").Value = 2.014742014742
$ws.Range("L17").Value = 16.957746478873
$ws.Range("M17").Value = 28.704277743335
$ws.Range("N17").Value = -49.909518639160
$ws.Range("C18").Value = 36
$ws.Range("D18").Value = 39
$ws.Range("E18").Value = -7.692307692307
$ws.Range("F18").Value = 138
$ws.Range("G18").Value = 152
$ws.Range("H18").Value = -9.210526315789
$ws.Range("I18").Value = 1983
$ws.Range("J18").Value = 2331
$ws.Range("K18").Value = -14.929214929214
$ws.Range("L18").Value = -4.983229516051
$ws.Range("M18").Value = -37.895396179141
$ws.Range("N18").Value = -83.331932419937
$ws.Range("C19").Value = 98
$ws.Range("D19").Value = 94
$ws.Range("E19").Value = 4.255319148936
$ws.Range("F19").Value = 429
$ws.Range("G19").Value = 414
$ws.Range("H19").Value = 3.623188405797
$ws.Range("I19").Value = 5649
$ws.Range("J19").Value = 5857
$ws.Range("K19").Value = -3.551306129417
$ws.Range("L19").Value = 17.982456140350
$ws.Range("M19").Value = 31.25
$ws.Range("N19").Value = -17.580974613364
$ws.Range("C20").Value = 32
$ws.Range("D20").Value = 25
$ws.Range("E20").Value = 28
$ws.Range("F20").Value = 156
$ws.Range("G20").Value = 153
$ws.Range("H20").Value = 1.960784313725
$ws.Range("I20").Value = 1839
$ws.Range("J20").Value = 1836
$ws.Range("K20").Value = 0.163398692810
$ws.Range("L20").Value = 21.788079470198
$ws.Range("M20").Value = 30.518097941802
$ws.Range("N20").Value = -80.498409331919
$ws.Range("C21").Value = 300
$ws.Range("D21").Value = 277
$ws.Range("E21").Value = 8.303249097472
$ws.Range("F21").Value = 1232
$ws.Range("G21").Value = 1174
$ws.Range("H21").Value = 4.940374787052
$ws.Range("I21").Value = 16391
$ws.Range("J21").Value = 16923
$ws.Range("K21").Value = -3.143650652957
$ws.Range("L21").Value = 14.462290502793
$ws.Range("M21").Value = 1.529980178394
$ws.Range("N21").Value = -69.710801071791
$ws.Range("C22").Value = 10
$ws.Range("D22").Value = 6
$ws.Range("E22").Value = 66.666666666666
$ws.Range("F22").Value = 30
$ws.Range("G22").Value = 27
$ws.Range("H22").Value = 11.111111111111
$ws.Range("I22").Value = 290
$ws.Range("J22").Value = 345
$ws.Range("K22").Value = -15.942028985507
$ws.Range("L22").Value = 0.346020761245
$ws.Range("M22").Value = -34.240362811791
$ws.Range("C23").Value = 24
$ws.Range("D23").Value = 27
$ws.Range("E23").Value = -11.111111111111
$ws.Range("F23").Value = 103
$ws.Range("G23").Value = 105
$ws.Range("H23").Value = -1.904761904761
$ws.Range("I23").Value = 1515
$ws.Range("J23").Value = 1486
$ws.Range("K23").Value = 1.951547779273
$ws.Range("L23").Value = 4.123711340206
$ws.Range("M23").Value = 30.490956072351
$ws.Range("C24").Value = 231
$ws.Range("D24").Value = 205
$ws.Range("E24").Value = 12.682926829268
$ws.Range("F24").Value = 980
$ws.Range("G24").Value = 971
$ws.Range("H24").Value = 0.926879505664
$ws.Range("I24").Value = 12052
$ws.Range("J24").Value = 13165
$ws.Range("K24").Value = -8.454234713254
$ws.Range("L24").Value = 15.164835164835
$ws.Range("M24").Value = 16.231073391841
$ws.Range("C25").Value = 111
$ws.Range("D25").Value = 93
$ws.Range("E25").Value = 19.354838709677
$ws.Range("F25").Value = 410
$ws.Range("G25").Value = 375
$ws.Range("H25").Value = 9.333333333333
$ws.Range("I25").Value = 5987
$ws.Range("J25").Value = 5765
$ws.Range("K25").Value = 3.850823937554
$ws.Range("L25").Value = 25.592615900985
$ws.Range("M25").Value = -22.377803708025
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = -55.555555555555
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = -3.703703703703
$ws.Range("I26").Value = 332
$ws.Range("J26").Value = 374
$ws.Range("K26").Value = -11.229946524064
$ws.Range("L26").Value = -10.512129380053
$ws.Range("C27").Value = 19
$ws.Range("D27").Value = 10
$ws.Range("E27").Value = 90
$ws.Range("F27").Value = 55
$ws.Range("G27").Value = 38
$ws.Range("H27").Value = 44.736842105263
$ws.Range("I27").Value = 619
$ws.Range("J27").Value = 596
$ws.Range("K27").Value = 3.859060402684
$ws.Range("L27").Value = -10.029069767441
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = -20
$ws.Range("F28").Value = 18
$ws.Range("G28").Value = 20
$ws.Range("H28").Value = -10
$ws.Range("I28").Value = 241
$ws.Range("J28").Value = 337
$ws.Range("K28").Value = -28.486646884273
$ws.Range("L28").Value = -40.493827160493
$ws.Range("M28").Value = -51.214574898785
$ws.Range("N28").Value = -86.801752464403
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = -40
$ws.Range("F29").Value = 15
$ws.Range("G29").Value = 19
$ws.Range("H29").Value = -21.052631578947
$ws.Range("I29").Value = 204
$ws.Range("J29").Value = 285
$ws.Range("K29").Value = -28.421052631578
$ws.Range("L29").Value = -38.738738738738
$ws.Range("M29").Value = -49.253731343283
$ws.Range("N29").Value = -87.591240875912
$ws.Range("C30").Value = 2
$ws.Range("F30").Value = 7
$ws.Range("H30").Value = 75
$ws.Range("I30").Value = 80
$ws.Range("K30").Value = -9.090909090909
$ws.Range("L30").Value = 33.333333333333
